# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For every data row (rows 2-24), a new leading error value is inserted
# into column B. All of the existing values in that row (previously in
# columns B..K) shift one column to the right (B->C, C->D, ... J->K),
# with whatever value was in the last occupied column of the row falling
# off the end (the row never grows past column K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to insert at column B for each row.
$newB = @{
    2  = 2.297389002388887
    3  = 8.826710628892494
    4  = -9.780318414391347
    5  = -1.200275438764269
    6  = 0.3719860057927588
    7  = -2.702915518772638
    8  = -0.2307826431404359
    9  = -0.5654386276933741
    10 = -0.6603092772102132
    11 = -0.15162438770796
    12 = -0.2053460154962278
    13 = 0.6162032393936197
    14 = 1.652643173475852
    15 = 0.3110387314724781
    16 = 0.2388379152847414
    17 = 0.6508000635779043
    18 = 0.2387740594105157
    19 = 0.3465902496671606
    20 = 0.00230005330798793
    21 = -0.1902738424076751
    22 = -0.3325070745318338
    23 = 0.1656141382254278
    24 = -0.09587373626955231
}

for ($row = 2; $row -le 24; $row++) {

    # Find the last occupied column in B:K for this row (K=11 .. B=2).
    $lastCol = 1
    for ($col = 11; $col -ge 2; $col--) {
        $v = $ws.Cells.Item($row, $col).Value2
        if ($v -ne $null) {
            $lastCol = $col
            break
        }
    }

    # Shift existing values one column to the right, starting from the
    # rightmost occupied column working back down to column C, so we
    # never clobber a value before it has been copied onward. Any value
    # that was already in column K is simply overwritten (dropped),
    # matching the row never extending past column K.
    for ($col = [Math]::Min($lastCol + 1, 11); $col -ge 3; $col--) {
        $srcVal = $ws.Cells.Item($row, $col - 1).Value2
        $ws.Cells.Item($row, $col).Value = $srcVal
    }

    # Insert the new leading value into column B.
    $ws.Cells.Item($row, 2).Value = $newB[$row]
}
